$d = $word.ActiveDocument

# Start from the end of the final existing paragraph ("This situation showcases...").
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)

# 1) Blank paragraph
$rng.InsertParagraphAfter()
$rng.Collapse(0)

# 2) "Assignment 3: Monitoring Natural Light for Energy Efficiency in a Smart Office Building"
#    (bold heading - text set now, bold formatting applied at the very end to avoid the
#    "sticky typing format" leaking into later paragraphs)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "Assignment 3: Monitoring Natural Light for Energy Efficiency in a Smart Office Building"
$headingRange = $headingPara.Range.Duplicate()

# 3) Blank paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)

# 4) Scenario paragraph
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "Scenario: You have been hired by a forward-thinking technology company to implement a smart lighting control system in their newly constructed office building. The company aims to maximize energy efficiency and create a comfortable workspace for employees."

# 5) Blank paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)

# 6) Situation Description:
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "Situation Description:"

# 7) Blank paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)

# 8) Large windows paragraph
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "The office building is equipped with large windows that allow ample natural light to enter."

# 9) Advanced lighting control system paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "To optimize energy consumption, the building uses an advanced lighting control system that adjusts the artificial lighting based on the available natural light."

# 10) ESP32 microcontrollers paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "ESP32 microcontrollers with LDR sensors are installed in each office space, near the windows."

# 11) LDR sensors paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "The LDR sensors continuously monitor the ambient light intensity, ranging from direct sunlight to cloudy conditions."

# 12) ESP32 devices send data paragraph
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "The ESP32 devices send the light intensity data to AWS IoT Core for real-time analysis and control."

# 13) Functionality:
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$d.Paragraphs.Last.Range.Text = "Functionality:"

# Apply bold (+ complex-script bold) to the Assignment 3 heading now that all the
# subsequent paragraphs already exist, so the formatting does not leak forward.
$headingRange.Bold = 1
$headingRange.BoldBi = 1

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
